$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the "Upload" row (row 8), shifting "Upload" and every row
# below it down by one. The new row becomes row 8, between "Save" (row 7) and the
# (now shifted) "Upload" row (row 9).
$ws.Rows.Item(8).Insert()

# Give the new row the same visual formatting (border/fill/font) as its sibling
# boolean-flag row "Save" directly above it.
$ws.Range("A7:G7").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row: label "Force" in column A, FALSE flags across B:G.
$ws.Range("A8").Value = "Force"
$ws.Range("B8:G8").Value = $false

# The sheet keeps its header block frozen. Since a row was inserted above the old
# freeze boundary (was between rows 9/10), move the freeze line down by one row
# (now between rows 10/11) so the same header block stays pinned.
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A11").Select()
$excel.ActiveWindow.FreezePanes = $true

# Leave the final selection on the newly inserted row's label cell.
$ws.Range("A9").Select()
